$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 164 (existing rows 164-179 shift down to 166-181)
$ws.Rows.Item(164).Insert()
$ws.Rows.Item(164).Insert()

# New row 164
$ws.Range("A164").Value = 3
$ws.Range("B164").Value = "Femacal de La Calera"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 44984
$ws.Range("E164").Value = 5
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100107
$ws.Range("H164").Value = "Otros"
$ws.Range("I164").Value = 100107011
$ws.Range("J164").Value = "Tuna"
$ws.Range("K164").Value = "Sin especificar"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 104
$ws.Range("N164").Value = 14000
$ws.Range("O164").Value = 15000
$ws.Range("P164").Value = 14481
$ws.Range("Q164").Value = "$/caja 18 kilos"
$ws.Range("R164").Value = "Provincia de Limarí"
$ws.Range("S164").Value = 804
$ws.Range("T164").Value = 18

# New row 165
$ws.Range("A165").Value = 3
$ws.Range("B165").Value = "Femacal de La Calera"
$ws.Range("C165").Value = "Coquimbo"
$ws.Range("D165").Value = 44984
$ws.Range("E165").Value = 5
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100107
$ws.Range("H165").Value = "Otros"
$ws.Range("I165").Value = 100107011
$ws.Range("J165").Value = "Tuna"
$ws.Range("K165").Value = "Sin especificar"
$ws.Range("L165").Value = "Segunda"
$ws.Range("M165").Value = 70
$ws.Range("N165").Value = 12000
$ws.Range("O165").Value = 12000
$ws.Range("P165").Value = 12000
$ws.Range("Q165").Value = "$/caja 18 kilos"
$ws.Range("R165").Value = "Provincia de Limarí"
$ws.Range("S165").Value = 667
$ws.Range("T165").Value = 18
